$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the J1 header cell value from "Studentoffer" to "Specialoffer"
$ws.Range("J1").Value = "Specialoffer"

# Update the sheet view: scroll so column B is the top-left visible column,
# and move the active selection to I8
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("I8").Select()
